$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "HUAMANI TORRES, LUIS RODRIGO - 06:30AM - 10:15AM"
$ws.Range("C1").Value = "MORENO CANCHANYA, ROSMERY - 10:30AM - 02:15PM"
$ws.Range("D1").Value = "YACILA GRANDEZ, RODRIGO ANDRE - 03:00PM - 06:45PM"
$ws.Range("E1").Value = "LA ROSA EUSEBIO, SHADIA SHAMIRA - 07:00PM - 10:45PM"
$ws.Range("B2").Value = "MARTINEZ PAZ, ROCIO ESPERANZA - 08:30AM - 05:30PM"
$ws.Range("C2").Value = "ZAVALA SOSA, NICOLE - 05:30PM - 09:15PM"
$ws.Range("B3").Value = "VEGA RIVAS, ANDREA FERNANDA - 09:30AM - 01:15PM"
$ws.Range("B4").Value = "AGUILAR SCHLAEFLI, STEPHANIE XIMENA - 07:00AM - 10:45AM"
$ws.Range("C4").Value = "Del Aguila Murayari, Darla - 10:45AM - 02:30PM"
$ws.Range("D4").Value = "HUAYNATES ALTAMIRANO, JIM HANS - 03:00PM - 06:45PM"
$ws.Range("E4").Value = "ARIAS MACHACUAY, SADELITH SORAGGI - 07:00PM - 10:45PM"
$ws.Range("B5").Value = "MENDOZA DIEGO, ZAIDA VANESSA - 08:45AM - 12:30PM"
$ws.Range("C5").Value = "BARRIENTOS JERI, MILAGROS NICOL - 01:00PM - 10:00PM"
$ws.Range("D5").Value = ""
$ws.Range("B6").Value = "MEZA MELO, NORMA FERNANDA - 08:45AM - 05:45PM"
$ws.Range("C6").Value = "MENDOZA CRUZ, LILIANA LILIANA - 05:45PM - 09:30PM"
$ws.Range("D6").Value = ""
$ws.Range("B7").Value = "YOVERA ROBLES, VICTOR EDUARDO - 09:00AM - 12:45PM"
$ws.Range("C7").Value = "SOTELO GONZALES, CAMILA SOFÍA - 01:00PM - 10:00PM"
$ws.Range("D7").Value = ""
$ws.Range("B8").Value = "AYALA MORA, CECILIA ROSARIO - 09:00AM - 12:45PM"
$ws.Range("C8").Value = "SAAVEDRA ALVAN, ANA MARIA - 01:00PM - 10:00PM"
$ws.Range("B9").Value = "VEGA CARDENAS, ANGELICA LOURDES - 08:00AM - 11:45AM"
$ws.Range("C9").Value = "ALVITE CORNEJO, ANGIE LUCERO - 12:00PM - 03:45PM"
$ws.Range("D9").Value = "VILCAPOMA CHILIN, JULISSA JAZMIN - 04:00PM - 07:45PM"
$ws.Range("B10").Value = "CHIARA LIMA, AUGUSTO SEBASTIAN - 08:00AM - 11:45AM"
$ws.Range("C10").Value = "FLORES PAREDES, LOURDES - 12:00PM - 02:45PM"
$ws.Range("D10").Value = "RIVERA CARREÑO, DIANA DESIRÉE - 03:30PM - 07:15PM"
$ws.Range("B11").Value = "MONTEZUMA DEJO, EVELYN BRUNELLA - 09:45AM - 01:30PM"
$ws.Range("C11").Value = "QUISPE MENDOZA, ANTONY MAURICIO - 02:00PM - 10:00PM"
$ws.Range("B12").Value = "MEZA PEREZ, JUAN CRISTOFER - 10:00AM - 01:45PM"
$ws.Range("C12").Value = "VILCHEZ CUBA, JACK ANTHONY - 02:00PM - 05:45PM"
$ws.Range("D12").Value = "CARDENAS RICAPA, FABRIZIO ESTEBAN - 06:00PM - 09:45PM"
$ws.Range("B13").Value = "HUAYANAY VELASCO, ATHINA - 10:00AM - 01:45PM"
$ws.Range("C13").Value = "PARICELA TINEO, JAIME DANIEL - 02:00PM - 05:45PM"
$ws.Range("D13").Value = "INGA DELGADO, CARLOS DANIEL - 06:00PM - 09:45PM"
$ws.Range("B14").Value = "CARHUARICRA ESPINOZA, FIORELLA NICOLL - 10:30AM - 02:15PM"
$ws.Range("C14").Value = "BRENIS LÁRTIGA, SEBASTIÁN - 04:00PM - 07:45PM"
$ws.Range("B15").Value = "CUSI QUISPE, ANDREA ESTEFANY - 10:30AM - 02:15PM"
$ws.Range("C15").Value = "YANAC DAVILA, GERALD RONNY - 05:00PM - 08:45PM"
